# Rotation Example.pptx - apply the commit's changes via PowerPoint COM interop
#
# Summary of changes:
#   1. Slide order: the slide that was at position 5 ("LookAt" content) and the
#      slide at position 6 ("네 번째 튜토리얼" content) swap places.
#   2. Slide 1 ("그림 8" picture, shape id 36): update the embedded picture's
#      description (original file path) and resize it slightly.
#   3. Slide 9 ("그림 7" picture, shape id 38): set a description (original file
#      path) that it didn't have before, nudge its left edge and resize it
#      slightly.
#
# NOTE on units: PowerPoint COM exposes Left/Top/Width/Height in points as a
# 32-bit float (Single), while the underlying OOXML stores EMU (1 pt = 12700
# EMU). Because of the float32 round-trip, plain "emu/12700.0" sometimes lands
# one EMU short after being truncated back down, so the literal point values
# below were chosen (by probing the runtime) to land exactly on the target EMU
# after that round trip.

$p = $ppt.ActivePresentation

# -----------------------------------------------------------------------
# 1) Swap the order of slides 5 and 6.
# -----------------------------------------------------------------------
$moved = $p.Slides.Item(5)
$moved.MoveTo(6)

# -----------------------------------------------------------------------
# 2) Slide 1 - picture "그림 8" (shape id 36): descr + size change.
# -----------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$pic36 = $slide1.Shapes.Item(5)
$pic36.AlternativeText = "C:/Users/nwond/AppData/Roaming/PolarisOffice/ETemp/20332_17372264/fImage12394608467.png"
$pic36.Width = 328.6500244140625
$pic36.Height = 104.55000305175781

# -----------------------------------------------------------------------
# 3) Slide 9 - picture "그림 7" (shape id 38): descr + position/size change.
# -----------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$pic38 = $slide9.Shapes.Item(5)
$pic38.AlternativeText = "C:/Users/nwond/AppData/Roaming/PolarisOffice/ETemp/19156_18797576/fImage421642018467.png"
$pic38.Left = 539.1500244140625
$pic38.Width = 324.20001220703125
$pic38.Height = 286.25
